$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Part 3"
$ws.Range("C4").Value = 10

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Part 4"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Part 5"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = (Get-Date -Year 1999 -Month 8 -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Part 6"
$ws.Range("C7").Value = 2

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Part 7"
$ws.Range("C8").Value = 3

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Part 8"
$ws.Range("C9").Value = 4

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Part 9"
$ws.Range("C10").Value = 5.5

$ws.Range("D7").Select()
